$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Update entry_made_on date format (column D) to yyyy/mm/dd, matching
#    the style used on the header (bold + vertical-center) for D1 and a
#    plain numeric date style for D2:D5.
# ---------------------------------------------------------------------
$ws.Range("D1").NumberFormat = "[$-14009]yyyy/mm/dd;@"
$ws.Range("D2:D5").NumberFormat = "[$-14009]yyyy/mm/dd;@"

# ---------------------------------------------------------------------
# 2. Update employee_code (B) and employee_name (C) for the 4 rows, plus
#    the one-off entry_made_on date and mobile_no on row 2.
# ---------------------------------------------------------------------
$ws.Range("B2").Value = 4000
$ws.Range("B3").Value = 4001
$ws.Range("B4").Value = 4002
$ws.Range("B5").Value = 4003

$ws.Range("C2").Value = "Satish Wakde"
$ws.Range("C3").Value = "Vinayak Mali"
$ws.Range("C4").Value = "Prathmesh Killedar"
$ws.Range("C5").Value = "Krishna Ware"

$ws.Range("D2").Value = 44523

$ws.Range("H2").Value = 918574968574

# ---------------------------------------------------------------------
# 3. Drop the pay_scale and payscale_per_hour columns. pay_scale sits at
#    AB, payscale_per_hour at AD; pay_scale_type (AC) shifts left into AB.
# ---------------------------------------------------------------------
$ws.Columns("AD:AD").Delete()
$ws.Columns("AB:AB").Delete()

# ---------------------------------------------------------------------
# 4. Update the view: drop the scrolled-away topLeftCell and move the
#    selection back to B5.
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("B5").Select()
